# dataset_nitrification.xlsx edit:
# - Rename Sheet1 -> "original", Sheet2 -> "idyno"
# - Update header labels on "original" with units
# - Populate "idyno" with a unit-converted copy (idynomics protocol units)
#   time (min) = original time (days) * 24 * 60
#   oxygen/ammonium/nitrite/nitrate (p.../fl) = original (mg/l) * 0.000001
# - Make "idyno" the active sheet/tab

$wb = $excel.ActiveWorkbook

$wsOriginal = $wb.Worksheets.Item(1)
$wsIdyno = $wb.Worksheets.Item(2)

$wsOriginal.Name = "original"
$wsIdyno.Name = "idyno"

# --- Update headers on "original" sheet (append units) ---
$wsOriginal.Range("A1").Value = "time (days)"
$wsOriginal.Range("B1").Value = "oxygen (mgO2/l)"
$wsOriginal.Range("C1").Value = "ammonium (mgN/l)"
$wsOriginal.Range("D1").Value = "nitrite (mgN/l)"
$wsOriginal.Range("E1").Value = "nitrate (mgN/l)"

# clear the old selection artifact and select the header row instead
[void]$wsOriginal.Range("A1:E1").Select()

# widen columns to fit the longer, unit-suffixed header labels (bestFit-style)
$wsOriginal.Range("B1").ColumnWidth = 15.166666666666666
$wsOriginal.Range("C1").ColumnWidth = 18.0
$wsOriginal.Range("D1").ColumnWidth = 13.333333333333334
$wsOriginal.Range("E1").ColumnWidth = 13.666666666666666

# --- Build "idyno" sheet: headers + formulas referencing "original" ---
$wsIdyno.Range("A1").Value = "time (min)"
$wsIdyno.Range("B1").Value = "oxygen (pgO2/fl)"
$wsIdyno.Range("C1").Value = "ammonium (pgN/fl)"
$wsIdyno.Range("D1").Value = "nitrite (pgN/fl)"
$wsIdyno.Range("E1").Value = "nitrate (pgN/fl)"

for ($r = 2; $r -le 36; $r++) {
    $wsIdyno.Range("A$r").Formula = "=original!A$r*24*60"
    $wsIdyno.Range("B$r").Formula = "=original!B$r*0.000001"
    $wsIdyno.Range("C$r").Formula = "=original!C$r*0.000001"
    $wsIdyno.Range("D$r").Formula = "=original!D$r*0.000001"
    $wsIdyno.Range("E$r").Formula = "=original!E$r*0.000001"
}

$wsIdyno.Range("A1").ColumnWidth = 10.166666666666666
$wsIdyno.Range("B1").ColumnWidth = 15.166666666666666
$wsIdyno.Range("C1").ColumnWidth = 18.0
$wsIdyno.Range("D1").ColumnWidth = 13.333333333333334
$wsIdyno.Range("E1").ColumnWidth = 13.666666666666666

[void]$wsIdyno.Range("E2").Select()
[void]$wsIdyno.Activate()
